# Updated cryptos list on Sat Jul 29 13:55:44 UTC 2023 with GitHub Actions
# Refresh price (D) / volume-change (E) figures and re-rank two tied pairs
# of rows (35/36 and 43/44) whose coins swapped order in the new ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.304.35"
$ws.Range("E2").Value = "  +0.06%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.78"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7096"
$ws.Range("E5").Value = "  -0.29%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.01"
$ws.Range("E6").Value = "  +0.08%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.00%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07793"
$ws.Range("E8").Value = "  +1.16%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3106"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.10"
$ws.Range("E10").Value = "  +1.45%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08418"
$ws.Range("E11").Value = "  +0.28%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.95"
$ws.Range("E12").Value = "  -0.31%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.239"
$ws.Range("E13").Value = "  +0.30%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7162"
$ws.Range("E14").Value = "  +0.50%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.14"
$ws.Range("E15").Value = "  -0.01%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.315.36"
$ws.Range("E16").Value = "  +0.11%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008298"
$ws.Range("E17").Value = "  +1.23%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.079"
$ws.Range("E18").Value = "  +2.35%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.58"
$ws.Range("E19").Value = "  -1.27%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.23"
$ws.Range("E20").Value = "  +0.68%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.117.39"
$ws.Range("E21").Value = "  -0.56%  "
# Row 22
$ws.Range("E22").Value = "  +0.02%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.756"
$ws.Range("E23").Value = "  -1.57%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1596"
$ws.Range("E25").Value = "  -1.69%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.49"
$ws.Range("E26").Value = "  -1.11%  "
# Row 27
$ws.Range("E27").Value = "  +0.26%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  -0.15%  "
# Row 29
$ws.Range("E29").Value = "  -0.25%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.407"
$ws.Range("E30").Value = "  +0.08%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.292"
$ws.Range("E31").Value = "  -1.16%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.324"
$ws.Range("E32").Value = "  +0.93%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05379"
$ws.Range("E33").Value = "  +3.97%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.950"
$ws.Range("E34").Value = "  +1.78%  "
# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7511"
$ws.Range("E35").Value = "  -2.99%  "
# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.178"
$ws.Range("E36").Value = "  +0.60%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  +0.03%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01885"
$ws.Range("E38").Value = "  +1.31%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.235.28"
$ws.Range("E39").Value = "  +6.90%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.729"
$ws.Range("E40").Value = "  +0.62%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.479"
$ws.Range("E41").Value = "  +1.22%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8935"
$ws.Range("E42").Value = "  +0.21%  "
# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.37"
$ws.Range("E43").Value = "  -1.18%  "
# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.95"
$ws.Range("E44").Value = "  +5.36%  "
# Row 45
$ws.Range("E45").Value = "  +0.03%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.022.40"
$ws.Range("E46").Value = "  +0.10%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000124"
$ws.Range("E47").Value = "  +7.54%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.799"
$ws.Range("E48").Value = "  -0.11%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5199"
$ws.Range("E49").Value = "  +0.15%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.442"
$ws.Range("E50").Value = "  +0.55%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4342"
$ws.Range("E51").Value = "  +1.01%  "
